$d = $word.ActiveDocument

$renames = @(
    @{ Old = "_Toc5728364"; New = "_Toc5728511" },
    @{ Old = "_Toc5728365"; New = "_Toc5728512" },
    @{ Old = "_Toc5728366"; New = "_Toc5728513" },
    @{ Old = "_Toc5728367"; New = "_Toc5728514" },
    @{ Old = "_Toc5728368"; New = "_Toc5728515" },
    @{ Old = "_Toc5728369"; New = "_Toc5728516" }
)

foreach ($item in $renames) {
    $bm = $d.Bookmarks.Item($item.Old)
    $rng = $bm.Range
    $bm.Delete()
    $d.Bookmarks.Add($item.New, $rng)
}
